$d = $word.ActiveDocument

$replacements = @(
    @{old="68×42=2856"; new="16×49=784"},
    @{old="63×96=6048"; new="47×28=1316"},
    @{old="75×14=1050"; new="13×24=312"},
    @{old="12×87=1044"; new="80×29=2320"},
    @{old="68×74=5032"; new="40×82=3280"},
    @{old="45×95=4275"; new="82×67=5494"},
    @{old="26×77=2002"; new="12×85=1020"},
    @{old="98×42=4116"; new="92×78=7176"},
    @{old="98×51=4998"; new="31×56=1736"},
    @{old="58×17=986";  new="43×85=3655"},
    @{old="95×29=2755"; new="18×85=1530"},
    @{old="86×17=1462"; new="85×66=5610"},
    @{old="90×13=1170"; new="17×52=884"},
    @{old="39×77=3003"; new="53×95=5035"},
    @{old="96×50=4800"; new="80×87=6960"},
    @{old="59×15=885";  new="66×93=6138"},
    @{old="39×67=2613"; new="86×35=3010"},
    @{old="67×96=6432"; new="29×70=2030"},
    @{old="51×72=3672"; new="21×30=630"},
    @{old="19×18=342";  new="55×41=2255"},
    @{old="82×62=5084"; new="27×74=1998"},
    @{old="86×31=2666"; new="20×19=380"},
    @{old="37×35=1295"; new="61×30=1830"},
    @{old="29×46=1334"; new="23×16=368"},
    @{old="46×72=3312"; new="87×98=8526"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
